$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("A1").Select() | Out-Null
